$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Admission No" column (D) for rows 2..67 with sequential
# admission numbers starting at 1001.
$admissionNo = 1001
for ($row = 2; $row -le 67; $row++) {
    $ws.Cells.Item($row, 4).Value = $admissionNo
    $admissionNo++
}

# Restore the active selection to K8.
$ws.Range("K8").Select()
